# Auto-update draw results: append the 2025-12-17 Pick 3 draw as a new last
# row (row 92) on the Results sheet, matching the pattern of every prior row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

$dateCell = $ws.Range("A" + $row)
$gameCell = $ws.Range("B" + $row)
$phaseCell = $ws.Range("C" + $row)
$resultCell = $ws.Range("D" + $row)
$insertedCell = $ws.Range("E" + $row)

# "2025-12-17" and "251217" look like a date / a plain number to Excel's
# input parser, so a plain .Value assignment would silently convert them to
# a date serial / a number instead of keeping them as literal text, like
# every other row in this sheet does. Force those two cells to Text first,
# write the values, then restore the original ("General"/default) style so
# the new row doesn't end up visibly formatted any differently than the
# rows above it.
$origStyle = $ws.Range("A" + ($row - 1)).Style

$dateCell.NumberFormat = "@"
$phaseCell.NumberFormat = "@"

$dateCell.Value = "2025-12-17"
$gameCell.Value = "Pick 3"
$phaseCell.Value = "251217"
$resultCell.Value = "6-8-7"
$insertedCell.Value = "2025-12-17T21:44:53.992+04:00"

$dateCell.Style = $origStyle
$phaseCell.Style = $origStyle
